# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
# Copy the formatting of the existing last header cell (AC1) onto the
# three new header cells so they match the rest of the header row
# (bold font, centered/top aligned, thin border).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# --- Data rows 2-59: team record is the same for every player row ---
$wins = 55
$losses = 106
$ties = 0

for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 30).Value2 = $wins    # column AD
    $ws.Cells.Item($r, 31).Value2 = $losses  # column AE
    $ws.Cells.Item($r, 32).Value2 = $ties    # column AF
}
